# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh).
# Coinranking snapshot values changed since the last run; rows 49/50 also
# swapped rank order (EnergySwap now above XinFinNetwork).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell -> new text value. Values are written as literal text (not
# numbers) so e.g. "12.50" keeps its trailing zero and "29.136.78" (a
# thousands-grouped price with two dots) isn't mis-parsed as a number.
$updates = [ordered]@{
    "D2" = "29.136.78"
    "E2" = "  -1.18%  "
    "D3" = "1.834.18"
    "E3" = "  -1.16%  "
    "D4" = "0.9988"
    "E4" = "  -0.02%  "
    "D5" = "240.03"
    "E5" = "  -1.91%  "
    "D6" = "0.6635"
    "E6" = "  -4.63%  "
    "E7" = "  +0.02%  "
    "D8" = "0.2954"
    "E8" = "  -3.75%  "
    "D9" = "0.07352"
    "E9" = "  -4.41%  "
    "E10" = "  -3.71%  "
    "D11" = "0.07678"
    "E11" = "  -1.28%  "
    "D12" = "1.839.58"
    "E12" = "  -1.00%  "
    "E13" = "  -2.69%  "
    "D14" = "0.6748"
    "E14" = "  -2.55%  "
    "D15" = "86.33"
    "E15" = "  -5.20%  "
    "D16" = "6.109"
    "E16" = "  -2.91%  "
    "D17" = "29.133.52"
    "E17" = "  -1.10%  "
    "D18" = "0.000008242"
    "E18" = "  -1.10%  "
    "D19" = "229.18"
    "E19" = "  -3.78%  "
    "D20" = "12.50"
    "E20" = "  -1.75%  "
    "D21" = "0.9995"
    "E21" = "  -0.03%  "
    "D22" = "7.297"
    "E22" = "  -4.16%  "
    "D23" = "0.9996"
    "E23" = "  -0.01%  "
    "D24" = "160.51"
    "E24" = "  +0.52%  "
    "D25" = "0.1419"
    "E25" = "  -5.16%  "
    "D26" = "8.662"
    "E26" = "  -2.55%  "
    "D27" = "18.02"
    "E28" = "  -1.86%  "
    "E29" = "  -0.19%  "
    "D30" = "4.098"
    "E30" = "  -1.17%  "
    "D31" = "1.196"
    "E31" = "  -0.18%  "
    "D32" = "0.05319"
    "E32" = "  +4.24%  "
    "D33" = "1.864"
    "E33" = "  -0.92%  "
    "E34" = "  -3.62%  "
    "E35" = "  -1.68%  "
    "D36" = "2.678"
    "E36" = "  -0.32%  "
    "D37" = "1.317.18"
    "E37" = "  +0.10%  "
    "D38" = "0.01802"
    "D39" = "2.710"
    "E39" = "  -0.29%  "
    "D40" = "0.9243"
    "E40" = "  -2.37%  "
    "D41" = "6.015"
    "E41" = "  +4.17%  "
    "D42" = "0.9984"
    "E42" = "  -0.13%  "
    "D43" = "103.34"
    "E43" = "  -2.60%  "
    "D44" = "1.982.51"
    "E44" = "  -0.87%  "
    "E46" = "  -3.01%  "
    "D47" = "1.759"
    "D48" = "63.51"
    "E48" = "  +0.81%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D49" = "9.278"
    "E49" = "  -5.13%  "
    "B50" = "XinFinNetwork"
    "C50" = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
    "D50" = "0.07489"
    "E50" = "  +11.13%  "
    "E51" = "  -0.12%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $newValue = $updates[$ref]

    # Preserve the cell's existing style/number-format. Forcing a
    # text numeric-format only for the instant of assignment stops Excel's
    # auto-detection from coercing a numeric-looking string (e.g. "240.03",
    # "12.50") into a real number (which would drop formatting like
    # trailing zeros). The original style is restored immediately after.
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = $originalStyle
}

Write-Host "Updated $($updates.Count) cells on $($ws.Name)"
